$wb = $excel.ActiveWorkbook

# --- "lines" sheet: insert a new "I_lim_A" column after v_nom_kv (before length_km) ---
# Shift the existing data (columns C..I) one column to the right by hand, instead of
# using Columns.Insert(), so the per-column <cols> width metadata (bound to columns
# D..H) stays attached to the same data and isn't shifted along with it.
$wsLines = $wb.Worksheets.Item("lines")
for ($r = 1; $r -le 4; $r++) {
    for ($c = 9; $c -ge 3; $c--) {
        $wsLines.Cells.Item($r, $c + 1).Value2 = $wsLines.Cells.Item($r, $c).Value2
    }
}

$wsLines.Range("C1").Value2 = "I_lim_A"
$wsLines.Range("C2").Value2 = 200
$wsLines.Range("C3").Value2 = 200
$wsLines.Range("C4").Value2 = 200

# --- "gens" sheet: drop the explicit number-format style on v_set_pu cells ---
$wsGens = $wb.Worksheets.Item("gens")
$wsGens.Range("C3").Style = "Normal"
$wsGens.Range("C4").Style = "Normal"

# --- Switch the active sheet/selection to "lines" (was "trafos") ---
$wsLines.Activate()
$wsLines.Range("D6").Select()

Write-Output "done"
